# Weekly update: insert two new price records (2023-12-20) at rows 35-36,
# pushing the existing rows 35-64 down to rows 37-66.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 35 (keeps formatting of the row above,
# matching Excel's default Insert behaviour).
$ws.Range("A35:A36").EntireRow.Insert()

# --- Row 35: Primera ---
$ws.Cells.Item(35, 1).Value = 7
$ws.Cells.Item(35, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(35, 3).Value = "Ñuble"
$ws.Cells.Item(35, 4).Value = [DateTime]"2023-12-20"
$ws.Cells.Item(35, 5).Value = 16
$ws.Cells.Item(35, 6).Value = "Fruta"
$ws.Cells.Item(35, 7).Value = 100101
$ws.Cells.Item(35, 8).Value = "Berries"
$ws.Cells.Item(35, 9).Value = 100101001
$ws.Cells.Item(35, 10).Value = "Arándano (blue)"
$ws.Cells.Item(35, 11).Value = "Sin especificar"
$ws.Cells.Item(35, 12).Value = "Primera"
$ws.Cells.Item(35, 13).Value = 100
$ws.Cells.Item(35, 14).Value = 4000
$ws.Cells.Item(35, 15).Value = 4000
$ws.Cells.Item(35, 16).Value = 4000
$ws.Cells.Item(35, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(35, 18).Value = "Región de Ñuble"
$ws.Cells.Item(35, 19).Value = 2000
$ws.Cells.Item(35, 20).Value = 2

# --- Row 36: Segunda ---
$ws.Cells.Item(36, 1).Value = 7
$ws.Cells.Item(36, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(36, 3).Value = "Ñuble"
$ws.Cells.Item(36, 4).Value = [DateTime]"2023-12-20"
$ws.Cells.Item(36, 5).Value = 16
$ws.Cells.Item(36, 6).Value = "Fruta"
$ws.Cells.Item(36, 7).Value = 100101
$ws.Cells.Item(36, 8).Value = "Berries"
$ws.Cells.Item(36, 9).Value = 100101001
$ws.Cells.Item(36, 10).Value = "Arándano (blue)"
$ws.Cells.Item(36, 11).Value = "Sin especificar"
$ws.Cells.Item(36, 12).Value = "Segunda"
$ws.Cells.Item(36, 13).Value = 100
$ws.Cells.Item(36, 14).Value = 3000
$ws.Cells.Item(36, 15).Value = 3000
$ws.Cells.Item(36, 16).Value = 3000
$ws.Cells.Item(36, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(36, 18).Value = "Región de Ñuble"
$ws.Cells.Item(36, 19).Value = 1500
$ws.Cells.Item(36, 20).Value = 2
